# Applies the tracked changes described by the commit diff to the document.
$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Highlight the run containing "МТС" in yellow.
#    Using Find/Replace with a Replacement.Highlight format avoids the
#    COM range-formatting bug that otherwise bleeds into neighbouring runs.
# ---------------------------------------------------------------------
$find1 = $d.Content
$find1.Find.ClearFormatting()
$find1.Find.Replacement.ClearFormatting()
$find1.Find.Text = "МТС"
$find1.Find.Replacement.Text = "МТС"
$find1.Find.Replacement.Highlight = 1
$find1.Find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Remove the old "_GoBack" bookmark that currently sits right after
#    "...е ниже 8.0" (before the ";" run). It gets re-created later, at
#    the very end of the document.
# ---------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 3) Split the sentence about mockup visual representations, dropping
#    "и функциональность" and re-wording the tail, as four separate runs.
# ---------------------------------------------------------------------
$find3 = $d.Content
$find3.Find.ClearFormatting()
$find3.Find.Replacement.ClearFormatting()
$oldSentence = ". Эти визуальные представления наглядно показывают структуру приложения, его основные элементы и функциональность."
$newSentence = ". Эти визуальные представления наглядно показывают структуру приложения"
$found3 = $find3.Find.Execute($oldSentence, $true, $false, $false, $false, $false, $true, 1, $false, $newSentence, 2)
if ($found3) {
    $find3.Collapse(0)
    $find3.InsertAfter(" ")
    $find3.Collapse(0)
    $find3.InsertAfter("и")
    $find3.Collapse(0)
    $find3.InsertAfter(" его основные элементы.")
}

# ---------------------------------------------------------------------
# 4) Tag the run that hosts the "Рисунок 2" deployment-diagram drawing
#    with an eastAsia language of ru-RU (matches the other figure run).
# ---------------------------------------------------------------------
$n = $d.Paragraphs.Count
for ($i = 1; $i -le $n; $i++) {
    $para = $d.Paragraphs.Item($i)
    $pr = $para.Range
    if ($pr.InlineShapes.Count -gt 0) {
        $capture = $pr.Duplicate
        $capture.MoveEnd(1, 3) | Out-Null
        $capture.Collapse(0)
        $checkFound = $capture.Find.Execute("Диаграмма развертывания", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if ($checkFound) {
            $pr.LanguageIDFarEast = "ru-RU"
        }
    }
}

# ---------------------------------------------------------------------
# 5) Final paragraph ("Требуется разработать БД для отдела ..."): add a
#    trailing "." run after the ERD sentence, then re-create "_GoBack"
#    as a collapsed bookmark right after it.
# ---------------------------------------------------------------------
$find5 = $d.Content
$find5.Find.ClearFormatting()
$find5.Find.Replacement.ClearFormatting()
$found5 = $find5.Find.Execute("показана физическая модель предметной области", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found5) {
    $find5.Collapse(0)
    $find5.InsertAfter(".")
    $find5.Collapse(0)
    # Insert then delete a throw-away character so the collapsed range used
    # for the bookmark never sits exactly on the (buggy) paragraph-end - 1
    # boundary while still ending up anchored right after the new ".".
    $find5.InsertAfter("Z")
    $find5.Collapse(1)
    if ($d.Bookmarks.Exists("_GoBack")) {
        $d.Bookmarks("_GoBack").Delete()
    }
    $d.Bookmarks.Add("_GoBack", $find5)
    $placeholderStart = $find5.Start
    $placeholder = $d.Range($placeholderStart, $placeholderStart + 1)
    $placeholder.Delete()
}

Write-Output "done"
